$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 148; this shifts rows 148:243 down to 149:244
$ws.Rows("148:148").Insert()

# Populate the new row 148 with the new data record
$ws.Range("A148").Value = 3
$ws.Range("B148").Value = "Femacal de La Calera"
$ws.Range("C148").Value = "Coquimbo"
$ws.Range("D148").Value = 44438
$ws.Range("E148").Value = 5
$ws.Range("F148").Value = 100112021
$ws.Range("G148").Value = "Ají"
$ws.Range("H148").Value = "Americana (o)"
$ws.Range("I148").Value = "Primera"
$ws.Range("J148").Value = 70
$ws.Range("K148").Value = 34000
$ws.Range("L148").Value = 35000
$ws.Range("M148").Value = 34500
$ws.Range("N148").Value = '$/caja 15 kilos'
$ws.Range("O148").Value = "Región de Arica y Parinacota"
$ws.Range("P148").Value = 2300
$ws.Range("Q148").Value = 15
$ws.Range("R148").Value = "Hortaliza"

# Make sure the date cell keeps the date/time number format used by the rest of column D
$ws.Range("D148").NumberFormat = $ws.Range("D149").NumberFormat
